$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$months = @("01","02","03","04","05","06","07","08","09","10","11","12")
$csaida = @(219,109,169,194,157,209,218,151,128,132,146,223)
$centrada = @(174,95,146,172,136,155,199,151,112,112,135,225)

$row = 44
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item($row, 1).Value = $months[$i]
    $ws.Cells.Item($row, 2).Value = "2018"
    $ws.Cells.Item($row, 3).Value = $csaida[$i]
    $ws.Cells.Item($row, 4).Value = $centrada[$i]
    $row = $row + 1
}

$ws.Range("D56").Select() | Out-Null
